$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D5").Value = 183.43
$ws.Range("D6").Value = 428.98
$ws.Range("D7").Value = 916.56
